$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$wsMeta = $wb.Worksheets.Item("Metadata")

# URL
$wsMeta.Range("B2").Value = "https://nphcda.gov.ng/immunizationIG/CodeSystem/nigeria-facility-type"

# Date
$wsMeta.Range("B8").Value = "2025-07-03T11:54:34+01:00"

# Count (keep as text, matching original shared-string type)
$wsMeta.Range("B23").Value = "'5"

# --- Concepts sheet updates ---
$wsConcepts = $wb.Worksheets.Item("Concepts")

# Row 2: clinic
$wsConcepts.Range("B2").Value = "clinic"
$wsConcepts.Range("C2").Value = "Clinic"
$wsConcepts.Range("D2").Value = "The Health facility is a primary health facility"

# Row 3: hospital
$wsConcepts.Range("B3").Value = "hospital"
$wsConcepts.Range("C3").Value = "Hospital"
$wsConcepts.Range("D3").Value = "The Health facility is a Secondary health facilities and General Hospitals"

# Row 4: health-post
$wsConcepts.Range("B4").Value = "health-post"
$wsConcepts.Range("C4").Value = "Health Post"
$wsConcepts.Range("D4").Value = "The Health facility is a primary health facility in category of dispensories"

# Row 5: chc (definition text unchanged)
$wsConcepts.Range("B5").Value = "chc"
$wsConcepts.Range("C5").Value = "Comprehensive Health Centers"

# Remove the old "public-general" and the two faith-based rows (rows 6 and 7, 8),
# leaving the existing "tertiary" row (old row 9) as the new row 6.
$wsConcepts.Rows("6:8").Delete()
